$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 236
$ws1.Range("F3").Value = 1057
$ws1.Range("F4").Value = 528
$ws1.Range("F5").Value = 13932
$ws1.Range("F7").Value = 562
$ws1.Range("F8").Value = 215
$ws1.Range("F9").Value = 1790
$ws1.Range("F10").Value = 174
$ws1.Range("F11").Value = 143
$ws1.Range("F12").Value = 94
$ws1.Range("F13").Value = 51
$ws1.Range("F14").Value = 537
$ws1.Range("F16").Value = 1
$ws1.Range("F18").Value = 13982
$ws1.Range("F20").Value = 630
$ws1.Range("F21").Value = 14986
$ws1.Range("F23").Value = 8288
$ws1.Range("F24").Value = 277
$ws1.Range("F26").Value = 27
$ws1.Range("F27").Value = 156
$ws1.Range("F28").Value = 426
$ws1.Range("F34").Value = 1040
$ws1.Range("F35").Value = 19
$ws1.Range("F38").Value = 406
$ws1.Range("F42").Value = 220
$ws1.Range("F43").Value = 393
$ws1.Range("F44").Value = 98
$ws1.Range("F45").Value = 5099
$ws1.Range("E14").Value = "2024.09.21 13:00-09.22 20:00"

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 236
$ws4.Range("F3").Value = 1057
$ws4.Range("F4").Value = 528
$ws4.Range("F5").Value = 13932
$ws4.Range("F7").Value = 562
$ws4.Range("F8").Value = 216
$ws4.Range("F9").Value = 1790
$ws4.Range("F10").Value = 174
$ws4.Range("F11").Value = 143
$ws4.Range("F12").Value = 94
$ws4.Range("F13").Value = 51
$ws4.Range("F14").Value = 537
$ws4.Range("F16").Value = 1
$ws4.Range("F18").Value = 13982
$ws4.Range("F20").Value = 630
$ws4.Range("F21").Value = 14986
$ws4.Range("F23").Value = 8288
$ws4.Range("F24").Value = 277
$ws4.Range("F26").Value = 27
$ws4.Range("F27").Value = 156
$ws4.Range("F28").Value = 426
$ws4.Range("F34").Value = 1040
$ws4.Range("F35").Value = 19
$ws4.Range("F40").Value = 406
$ws4.Range("F44").Value = 220
$ws4.Range("F45").Value = 393
$ws4.Range("F46").Value = 98
$ws4.Range("F47").Value = 5099
$ws4.Range("E14").Value = "2024.09.21 13:00-09.22 20:00"
